$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) values stay stored as text, matching the
# original inline-string cells, rather than being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Updated cryptos list values (Price / Volume(1h)) per the latest scrape
$ws.Range("D2").Value = '68.157.78'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '3.274.22'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '587.53'
$ws.Range("E5").Value = '  +1.97%  '
$ws.Range("D6").Value = '185.54'
$ws.Range("E6").Value = '  +3.92%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +4.78%  '
$ws.Range("D10").Value = '6.74'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("D12").Value = '3.842.60'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '28.63'
$ws.Range("E14").Value = '  +2.69%  '
$ws.Range("D15").Value = '68.142.78'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = '3.270.32'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '13.66'
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").Value = '382.61'
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '71.43'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("E24").Value = '  +2.78%  '
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").Value = '0.192'
$ws.Range("E26").Value = '  +7.00%  '
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("E29").Value = '  +4.43%  '
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").Value = '22.96'
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("D32").Value = '7.21'
$ws.Range("E32").Value = '  +6.00%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D36").Value = '163.45'
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("D39").Value = '6.77'
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").Value = '26.65'
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("D41").Value = '2.66'
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("D42").Value = '4.61'
$ws.Range("E42").Value = '  +4.90%  '
$ws.Range("D43").Value = '25.67'
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("E45").Value = '  +2.78%  '
$ws.Range("D46").Value = '2.637.95'
$ws.Range("E46").Value = '  -4.38%  '
$ws.Range("D47").Value = '342.47'
$ws.Range("E47").Value = '  -3.00%  '
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("D49").Value = '32.34'
$ws.Range("E49").Value = '  +5.29%  '
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("E51").Value = '  -0.17%  '

# Restore default (unstyled) formatting on the Price column so only the
# cell contents change, consistent with the rest of the sheet.
$ws.Range("D2:D51").Style = "Normal"
